# Fix Groupings for 2017
# - Change oversights label
# - Change t3* question labels to numbers so they work in the report generator
# - Make "QuestionsGroups 2017" the active/selected sheet (with C12 selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuestionsGroups 2017")

# Update the oversight institutions group label (B6)
$ws.Range("B6").Value2 = "Role and Effectiveness of Oversight Institutions"

# Replace the old "t3*" style question labels with numeric ones
$ws.Range("C12").Value2 = "1-53, 144"
$ws.Range("C13").Value2 = "59-63, 145"
$ws.Range("C15").Value2 = "68-75, 146"
$ws.Range("C16").Value2 = "76-83, 147"
$ws.Range("C17").Value2 = "84-96, 148"
$ws.Range("C18").Value2 = "97-102, 149"
$ws.Range("C11").Value2 = "54-58, 143"

# Make "QuestionsGroups 2017" sheet the active tab, with C12 selected
$ws.Activate()
$ws.Range("C12").Select()
